$p = $ppt.ActivePresentation

for ($i = 1; $i -le 4; $i++) {
    $s = $p.Slides.Item($i)
    $s.Shapes.Item(1).TextFrame.TextRange.Text = "Slide $i"
}
